# Add a blank line, another blank line, then a new command line
# ("webpack -c webpack.config.client.js") right after the existing
# "npm install react-router react-router-dom -save" paragraph, matching
# the formatting (Arial, 10.5pt, black paragraph shading) already used
# by the other command-line paragraphs in that block.

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find so the script is resilient to
# exact paragraph-index drift.
$rng = $d.Content
$found = $rng.Find.Execute("npm install react-router react-router-dom", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'npm install react-router' paragraph"
}

# Expand the found hit to cover the whole paragraph (this also picks up
# the paragraph's pPr/rPr formatting that InsertParagraphAfter copies
# forward into the new paragraphs).
[void]$rng.Expand(4)  # wdParagraph

# Remember the 1-based index of this paragraph so we can address the
# newly-created ones afterwards.
$preceding = $d.Range(0, $rng.Start)
$anchorIndex = $preceding.Paragraphs.Count + 1

# Insert three new paragraphs right after the anchor paragraph. Each
# call inserts immediately after the anchor (pushing the previous
# insertion further down), so after three calls the order below the
# anchor paragraph is: blank, blank, webpack-command.
[void]$rng.InsertParagraphAfter()
[void]$rng.InsertParagraphAfter()
[void]$rng.InsertParagraphAfter()

$webpackPara = $d.Paragraphs.Item($anchorIndex + 3)
$webpackPara.Range.Text = "webpack -c webpack.config.client.js"
